$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(130).Insert()

$ws.Cells.Item(130, 1).Value = 4
$ws.Cells.Item(130, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value = "Los Lagos"
$ws.Cells.Item(130, 4).Value = 45009
$ws.Cells.Item(130, 5).Value = 10
$ws.Cells.Item(130, 6).Value = 100112052
$ws.Cells.Item(130, 7).Value = "Albahaca"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 100
$ws.Cells.Item(130, 11).Value = 6500
$ws.Cells.Item(130, 12).Value = 8000
$ws.Cells.Item(130, 13).Value = 7250
$ws.Cells.Item(130, 14).Value = "`$/docena de matas"
$ws.Cells.Item(130, 15).Value = "Región Metropolitana"
$ws.Cells.Item(130, 16).Value = 1208
$ws.Cells.Item(130, 17).Value = 6
$ws.Cells.Item(130, 18).Value = "Hortaliza"
